# Applies the "Add files via upload" edit to the "Phase I" sheet:
#  - Relabel the G1 header from "Priority" to "Priority (1 is highest 3 is lowest)"
#  - Fill in the Assigned-To (F) and Priority (G) columns for several requirement rows
#  - Fix the fill colour on D99 so it matches its "complete" text (red -> green)
#  - Widen column G and move the saved view/selection to around G104

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Phase I")

# Assigned To / Priority values for the newly triaged requirement rows
# ("Ethan" is introduced first so it lands before the new header text in the
# shared-string table, matching the authored workbook.)
$ws.Range("G10").Value = 3

$ws.Range("G14").Value = 3

$ws.Range("F21").Value = "Ethan"
$ws.Range("G21").Value = 3

$ws.Range("F22").Value = "Ben"
$ws.Range("G22").Value = 3

$ws.Range("F24").Value = "Ethan"
$ws.Range("G24").Value = 3

$ws.Range("F25").Value = "Ben"
$ws.Range("G25").Value = 3

$ws.Range("F27").Value = "Ethan"
$ws.Range("G27").Value = 1

$ws.Range("F28").Value = "Ben"
$ws.Range("G28").Value = 1

$ws.Range("F30").Value = "Ben"
$ws.Range("G30").Value = 1

$ws.Range("F32").Value = "Ethan"
$ws.Range("G32").Value = 2

$ws.Range("F39").Value = "Ben"
$ws.Range("G39").Value = 3

$ws.Range("F40").Value = "Ben"
$ws.Range("G40").Value = 3

$ws.Range("F41").Value = "Ben"
$ws.Range("G41").Value = 3

$ws.Range("F42").Value = "Ben"
$ws.Range("G42").Value = 3

$ws.Range("F101").Value = "Ethan"
$ws.Range("G101").Value = 2

$ws.Range("F104").Value = "Ben"
$ws.Range("G104").Value = 2

# Header text update (this also retires the old standalone "Priority" shared string)
$ws.Range("G1").Value = "Priority (1 is highest 3 is lowest)"

# D99 is labelled "complete" but was still shaded red (incomplete) - recolour it green
# to match D3/D4/etc, by copying the known-good green style from D3.
$ws.Range("D3").Copy()
$ws.Range("D99").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Widen the new Priority column and restore the author's saved scroll/selection position
$ws.Columns.Item(7).ColumnWidth = 30.666666666666668

$ws.Application.GoTo($ws.Range("G104"), $true)
$ws.Range("G104").Select() | Out-Null
